$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = "Normal"
}

$ws.Range("D2").Value = "42.768.84"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "2.295.47"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  +0.23%  "
Set-TextValue $ws.Range("D5") "300.42"
$ws.Range("E5").Value = "  -2.64%  "
Set-TextValue $ws.Range("D6") "98.16"
$ws.Range("E6").Value = "  -7.27%  "
Set-TextValue $ws.Range("D7") "0.501"
$ws.Range("E7").Value = "  -4.71%  "
$ws.Range("E8").Value = "  +0.13%  "
Set-TextValue $ws.Range("D9") "0.497"
$ws.Range("E9").Value = "  -3.97%  "
Set-TextValue $ws.Range("D10") "34.43"
$ws.Range("E10").Value = "  -4.72%  "
Set-TextValue $ws.Range("D11") "0.0784"
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.663.94"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "6.63"
$ws.Range("E14").Value = "  -4.95%  "
Set-TextValue $ws.Range("D15") "15.34"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "2.340.12"
$ws.Range("E16").Value = "  +1.34%  "
Set-TextValue $ws.Range("D17") "0.786"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").Value = "42.738.94"
$ws.Range("E18").Value = "  -1.54%  "
Set-TextValue $ws.Range("D19") "11.51"
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("D20").Value = "0.0₃0893"
$ws.Range("E20").Value = "  -3.22%  "
Set-TextValue $ws.Range("D21") "5.96"
$ws.Range("E21").Value = "  -3.74%  "
Set-TextValue $ws.Range("D22") "67.35"
$ws.Range("E22").Value = "  -0.97%  "
Set-TextValue $ws.Range("D23") "234.41"
$ws.Range("E23").Value = "  -2.78%  "
Set-TextValue $ws.Range("D24") "1.93"
$ws.Range("E24").Value = "  -4.95%  "
Set-TextValue $ws.Range("D25") "2.49"
$ws.Range("E25").Value = "  -4.47%  "
$ws.Range("E26").Value = "  -0.05%  "
Set-TextValue $ws.Range("D27") "24.46"
$ws.Range("E27").Value = "  -2.42%  "
Set-TextValue $ws.Range("D28") "2.16"
$ws.Range("E28").Value = "  -2.25%  "
Set-TextValue $ws.Range("D29") "33.88"
$ws.Range("E29").Value = "  -7.04%  "
Set-TextValue $ws.Range("D30") "163.17"
$ws.Range("E30").Value = "  +0.05%  "
Set-TextValue $ws.Range("D31") "9.03"
$ws.Range("E31").Value = "  -5.83%  "
Set-TextValue $ws.Range("D32") "0.999"
$ws.Range("E32").Value = "  -0.02%  "
Set-TextValue $ws.Range("D33") "4.95"
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("E34").Value = "  -4.67%  "
Set-TextValue $ws.Range("D35") "4.41"
$ws.Range("E35").Value = "  -4.13%  "
Set-TextValue $ws.Range("D36") "16.47"
$ws.Range("E36").Value = "  -9.81%  "
Set-TextValue $ws.Range("D37") "0.0689"
$ws.Range("E37").Value = "  -6.11%  "
Set-TextValue $ws.Range("D38") "2.85"
$ws.Range("E38").Value = "  -5.38%  "
Set-TextValue $ws.Range("D39") "1.77"
$ws.Range("E39").Value = "  -4.89%  "
Set-TextValue $ws.Range("D40") "0.0997"
$ws.Range("E40").Value = "  -5.15%  "
Set-TextValue $ws.Range("D41") "0.109"
$ws.Range("E41").Value = "  -5.02%  "
Set-TextValue $ws.Range("D42") "2.51"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").Value = "1.951.35"
$ws.Range("E43").Value = "  -0.59%  "
Set-TextValue $ws.Range("D44") "0.0278"
$ws.Range("E44").Value = "  -3.84%  "
Set-TextValue $ws.Range("D45") "18.28"
$ws.Range("E45").Value = "  -3.05%  "
Set-TextValue $ws.Range("D46") "10.09"
$ws.Range("E46").Value = "  -1.91%  "
Set-TextValue $ws.Range("D47") "2.83"
$ws.Range("E47").Value = "  -7.21%  "
Set-TextValue $ws.Range("D48") "54.46"
$ws.Range("E48").Value = "  -6.33%  "
$ws.Range("D49").Value = "2.529.65"
$ws.Range("E49").Value = "  -0.21%  "
Set-TextValue $ws.Range("D50") "2.81"
$ws.Range("E50").Value = "  -4.84%  "
Set-TextValue $ws.Range("D51") "4.65"
$ws.Range("E51").Value = "  -1.53%  "
